$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1:XFD1048576").Select() | Out-Null

$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "summerToursTest"

$ws2.Hyperlinks.Delete()
$ws1.Range("B2").Copy($ws2.Range("B2"))
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:amanullah.a@gmail.com")

$ws2.Range("A2").Value = "Amanullah Akbar Ali"
$ws2.Range("B2").Value = "amanullah.a@gmail.com"
$ws2.Range("C2").Value = "9943357865"
$ws2.Range("D2").Value = "Hi This is Amanullah from Erode"

$ws2.Range("C12").Select() | Out-Null
